$d = $word.ActiveDocument

# 1) Fine Amount table cell: "$ 50" -> "$ 0"
foreach ($tbl in $d.Tables) {
    foreach ($row in $tbl.Rows) {
        $labelCell = $row.Cells.Item(1)
        $labelText = $labelCell.Range.Text
        if ($labelText -match "Fine Amount") {
            $valueCell = $row.Cells.Item(2)
            $valueCell.Range.Find.Execute("`$ 50", $true, $false, $false, $false, $false,
                                           $true, 1, $false, "`$ 0", 2)
        }
    }
}

# 2) Remove "Defendant's report date is January 01, 2000, at 08:30 AM." sentence entirely.
$d.Content.Find.Execute("Defendant" + [char]0x2019 + "s report date is January 01, 2000, at 08:30 AM.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2)

# 3) "The jail days imposed shall be served as consecutive days." -> "... as None. None."
$d.Content.Find.Execute(" consecutive days. ", $true, $false, $false, $false, $false,
                         $true, 1, $false, " None. None.", 2)
